$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting of the last existing data row down into the new row
$ws.Range("A55:E55").Copy()
$ws.Range("A56:E56").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Fill in the new row of data (2020-05-06)
$ws.Cells.Item(56, 1).Value = 43957
$ws.Cells.Item(56, 2).Value = 30303
$ws.Cells.Item(56, 3).Value = 2253
$ws.Cells.Item(56, 4).Value = 64
$ws.Cells.Item(56, 5).Value = 4917

# Grow the table so the new row is included
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E56"))

# Update view state to match
[void]$ws.Range("E56").Select()
